$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated data for rows 2-12 (columns A, C, D). Column B and row1 headers are unchanged.
$data = @(
  @{ Row=2;  A=0.1162216402902751; C=1.143548406185241; D=0.8282530991886153; E="/" },
  @{ Row=3;  A=0.1498031341417626; C=1.159955780764616; D=0.8057748446783823; E="/" },
  @{ Row=4;  A=0.296001417174004;  C=1.231821262345502; D=0.7193136704834024; E="/" },
  @{ Row=5;  A=0.4402667545719619; C=1.30344730518867;  D=0.6500612503083237; E="/" },
  @{ Row=6;  A=0.5592346242831165; C=1.363607283954145; D=0.6016703844379322; E="/" },
  @{ Row=7;  A=0.6846514270039978; C=1.427045742332629; D=0.5594648600787343; E=" " },
  @{ Row=8;  A=0.6846514270039978; C=1.412696676147325; D=0.5826988512254297; E=" " },
  @{ Row=9;  A=0.6846514270039978; C=1.399674171024946; D=0.6056513339984361; E=" " },
  @{ Row=10; A=0.6846514270039978; C=1.387771373107714; D=0.628362888747945;  E=" " },
  @{ Row=11; A=0.6846514270039978; C=1.37682585693647;  D=0.6508647199596943; E=" " },
  @{ Row=12; A=0.6846514270039978; C=1.366707751087719; D=0.6731811812813191; E=" " }
)

foreach ($rowData in $data) {
  $r = $rowData.Row
  $ws.Cells.Item($r, 1).Value = $rowData.A
  $ws.Cells.Item($r, 3).Value = $rowData.C
  $ws.Cells.Item($r, 4).Value = $rowData.D
  $ws.Cells.Item($r, 5).Value = $rowData.E
}

# Remove the now-obsolete rows 13-16 (original table had 16 rows, new table has 12).
$ws.Rows("13:16").Delete()
